# Updated symbol list on Tue Dec 13 11:32:21 UTC 2022 with GitHub Actions
#
# Refreshed crypto price snapshot: column D ("Price") values are stored as
# text in this sheet, so each price cell is pre-formatted as Text ("@")
# before the new value is written — this keeps the cell a text value
# instead of letting Excel auto-coerce the numeric-looking string into a
# Number. Column E ("Volume(1h)") values are plain labels that sometimes
# carry a "Worstin24h" suffix flagging the day's worst performer; that
# flag moved off UpBots (row 28) and on to CEJI (row 43) in this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2"  "268.43"
Set-TextValue "D3"  "21.56"
Set-TextValue "D4"  "6.208"
Set-TextValue "D5"  "0.06172"
Set-TextValue "D6"  "3.564"
Set-TextValue "D7"  "6.552"
Set-TextValue "D8"  "1.366"
Set-TextValue "D9"  "0.8238"
Set-TextValue "D10" "0.01347"
Set-TextValue "D11" "0.1560"
Set-TextValue "D12" "0.08153"
Set-TextValue "D13" "0.03314"
Set-TextValue "D14" "0.03178"
Set-TextValue "D15" "0.09272"
Set-TextValue "D16" "3.749"
Set-TextValue "D17" "0.001622"
Set-TextValue "D18" "0.04670"
Set-TextValue "D19" "0.006406"
Set-TextValue "D21" "0.001067"
Set-TextValue "D23" "3.737"
Set-TextValue "D24" "2.321"
Set-TextValue "D25" "0.3303"

Set-TextValue "D28" "0.0001616"
$ws.Range("E28").Value = "27UpBotsUBXT"

Set-TextValue "D40" "0.04658"
Set-TextValue "D41" "0.006966"
Set-TextValue "D42" "0.1132"

Set-TextValue "D43" "0.003446"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Set-TextValue "D44" "0.01183"
Set-TextValue "D45" "0.00006074"
Set-TextValue "D46" "0.0009884"
Set-TextValue "D48" "0.7809"
Set-TextValue "D49" "0.002439"
